# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> bound to the slide master ("Integral" palette)
#   ppt/theme/theme2.xml -> bound to the notes master  ("Office Theme" palette)
#
# The target revision swaps the two themes' colour content: the slide
# master's theme (theme1.xml) becomes the "Office Theme" palette, while the
# notes master's theme (theme2.xml) becomes the "Integral" palette. Font
# scheme / format scheme are already byte-identical between the two themes,
# so the only observable difference is the 12-slot colour scheme (plus the
# cosmetic <a:theme>/<a:clrScheme> name attributes).
#
# This COM host only exposes one editable theme-colour surface -
# Slide.ThemeColorScheme - and it always targets the slide master's theme
# part (theme1.xml); there is no reachable object for the notes master's
# separate theme part. So we push the "Office Theme" palette (which
# theme1.xml's sibling theme2.xml currently holds) onto theme1.xml via that
# API, slot by slot, using the documented dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink -> Item(1..12) ordering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RgbVal([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

# Target palette = the "Office Theme" colours (currently in theme2.xml).
$tcs.Item(1).RGB  = RgbVal 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RgbVal 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RgbVal 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RgbVal 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RgbVal 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RgbVal 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RgbVal 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RgbVal 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RgbVal 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RgbVal 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RgbVal 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RgbVal 0x95 0x4F 0x72   # folHlink
